$d = $word.ActiveDocument

$replacements = @(
    @("233×2=466", "558×3=1674"),
    @("436×3=1308", "788×7=5516"),
    @("921×9=8289", "733×5=3665"),
    @("207×4=828", "152×4=608"),
    @("988×7=6916", "377×6=2262"),
    @("150×5=750", "973×6=5838"),
    @("638×6=3828", "201×7=1407"),
    @("126×8=1008", "483×9=4347"),
    @("858×5=4290", "512×9=4608"),
    @("462×9=4158", "479×7=3353"),
    @("341×4=1364", "101×3=303"),
    @("103×6=618", "157×5=785"),
    @("868×5=4340", "671×6=4026"),
    @("301×4=1204", "612×4=2448"),
    @("233×8=1864", "799×3=2397"),
    @("293×9=2637", "231×8=1848"),
    @("758×5=3790", "540×2=1080"),
    @("244×6=1464", "892×7=6244"),
    @("179×5=895", "720×2=1440"),
    @("152×6=912", "867×3=2601"),
    @("616×9=5544", "177×8=1416"),
    @("720×7=5040", "233×5=1165"),
    @("621×7=4347", "293×2=586"),
    @("817×2=1634", "903×3=2709"),
    @("198×4=792", "766×6=4596")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
